# Regenerate save_data to use K (strikeouts) instead of Strike# and
# recompute the std/mean derived "K" column (column G) values written
# from the newly calculated s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column G ("K")
$kValues = @{
    2  = 0
    3  = 1
    5  = 2
    6  = 0
    7  = 2
    8  = 1
    9  = 1
    10 = 0
    11 = 0
    12 = 0
    15 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
